# B1--and-B2-PowerPoint.pptx edit
#
# The authoritative diff changes the <a:tableStyleId> of the table on
# slide 5 (the "B1- TYPES OF FINANCIAL DOCUMENTS" slide) from the
# built-in "No Style, Table Grid" GUID to "Light Style 1 - Accent 1"
# ({AE711E26-853F-4949-9850-5BA52A44F501}).
#
# PowerPoint's Table object exposes the applied table style through the
# read-only `Table.Style` property; changing it is done with
# `Table.ApplyStyle("{GUID}")`.

$p = $ppt.ActivePresentation

$targetStyleId = "{AE711E26-853F-4949-9850-5BA52A44F501}"

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if ($shape.HasTable) {
            $table = $shape.Table
            $table.ApplyStyle($targetStyleId)
        }
    }
}
